$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97   # AD
    $ws.Cells.Item($r, 31).Value = 66   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
